$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4,2).Value = 1148439
$ws.Cells.Item(4,3).Value = 17409
$ws.Cells.Item(4,4).Value = 162114
$ws.Cells.Item(4,5).Value = 919521
$ws.Cells.Item(4,7).Value = 1051
$ws.Cells.Item(4,8).Value = 66804
$ws.Cells.Item(8,2).Value = 168396
$ws.Cells.Item(8,3).Value = 1050
$ws.Cells.Item(8,4).Value = 50562
$ws.Cells.Item(8,5).Value = 93074
$ws.Cells.Item(8,6).Value = 3827
$ws.Cells.Item(8,7).Value = 166
$ws.Cells.Item(8,8).Value = 24760
$ws.Cells.Item(13,2).Value = 92865
$ws.Cells.Item(13,3).Value = 756
$ws.Cells.Item(13,5).Value = 48335
$ws.Cells.Item(13,7).Value = 81
$ws.Cells.Item(13,8).Value = 6491
$ws.Cells.Item(27,2).Value = 19022
$ws.Cells.Item(27,3).Value = 930
$ws.Cells.Item(27,5).Value = 13832
$ws.Cells.Item(27,7).Value = 20
$ws.Cells.Item(27,8).Value = 437
$ws.Cells.Item(51,1).Value = 'Sudafrica'
$ws.Cells.Item(51,2).Value = 6336
$ws.Cells.Item(51,3).Value = 385
$ws.Cells.Item(51,4).Value = 2549
$ws.Cells.Item(51,5).Value = 3664
$ws.Cells.Item(51,6).Value = 36
$ws.Cells.Item(51,7).Value = 7
$ws.Cells.Item(51,8).Value = 123
$ws.Cells.Item(52,1).Value = 'Egipto'
$ws.Cells.Item(52,2).Value = 6193
$ws.Cells.Item(52,3).Value = 298
$ws.Cells.Item(52,4).Value = 1522
$ws.Cells.Item(52,5).Value = 4256
$ws.Cells.Item(52,6).Value = 0
$ws.Cells.Item(52,7).Value = 9
$ws.Cells.Item(52,8).Value = 415
$ws.Cells.Item(53,1).Value = 'Malasia'
$ws.Cells.Item(53,2).Value = 6176
$ws.Cells.Item(53,3).Value = 105
$ws.Cells.Item(53,4).Value = 4326
$ws.Cells.Item(53,5).Value = 1747
$ws.Cells.Item(53,6).Value = 31
$ws.Cells.Item(53,8).Value = 103
$ws.Cells.Item(79,1).Value = 'Cuba'
$ws.Cells.Item(79,2).Value = 1611
$ws.Cells.Item(79,3).Value = 74
$ws.Cells.Item(79,4).Value = 765
$ws.Cells.Item(79,5).Value = 780
$ws.Cells.Item(79,6).Value = 10
$ws.Cells.Item(79,7).Value = 2
$ws.Cells.Item(79,8).Value = 66
$ws.Cells.Item(80,1).Value = 'Bulgaria'
$ws.Cells.Item(80,2).Value = 1594
$ws.Cells.Item(80,3).Value = 39
$ws.Cells.Item(80,4).Value = 287
$ws.Cells.Item(80,5).Value = 1235
$ws.Cells.Item(80,6).Value = 43
$ws.Cells.Item(80,7).Value = 4
$ws.Cells.Item(80,8).Value = 72
$ws.Cells.Item(81,1).Value = 'Guinea'
$ws.Cells.Item(81,4).Value = 342
$ws.Cells.Item(81,5).Value = 1188
$ws.Cells.Item(81,6).Value = 0
$ws.Cells.Item(81,8).Value = 7
$ws.Cells.Item(98,2).Value = 747
$ws.Cells.Item(98,3).Value = 2
$ws.Cells.Item(98,4).Value = 472
$ws.Cells.Item(98,5).Value = 231
$ws.Cells.Item(98,7).Value = 1
$ws.Cells.Item(98,8).Value = 44
$ws.Cells.Item(100,1).Value = 'Costa Rica'
$ws.Cells.Item(100,2).Value = 733
$ws.Cells.Item(100,3).Value = 8
$ws.Cells.Item(100,4).Value = 372
$ws.Cells.Item(100,5).Value = 355
$ws.Cells.Item(100,6).Value = 6
$ws.Cells.Item(100,8).Value = 6
$ws.Cells.Item(101,1).Value = 'Niger'
$ws.Cells.Item(101,2).Value = 728
$ws.Cells.Item(101,4).Value = 478
$ws.Cells.Item(101,5).Value = 217
$ws.Cells.Item(101,6).Value = 0
$ws.Cells.Item(101,8).Value = 33
$ws.Cells.Item(118,4).Value = 141
$ws.Cells.Item(118,5).Value = 295
$ws.Cells.Item(142,4).Value = 132
$ws.Cells.Item(142,5).Value = 12
$ws.Cells.Item(150,1).Value = 'Republica del Chad'
$ws.Cells.Item(150,2).Value = 117
$ws.Cells.Item(150,3).Value = 44
$ws.Cells.Item(150,4).Value = 39
$ws.Cells.Item(150,5).Value = 68
$ws.Cells.Item(150,7).Value = 5
$ws.Cells.Item(150,8).Value = 10
$ws.Cells.Item(151,1).Value = 'Trinidad yTobago'
$ws.Cells.Item(151,2).Value = 116
$ws.Cells.Item(151,4).Value = 87
$ws.Cells.Item(151,5).Value = 21
$ws.Cells.Item(151,6).Value = 0
$ws.Cells.Item(151,8).Value = 8
$ws.Cells.Item(152,1).Value = 'Bermudas'
$ws.Cells.Item(152,2).Value = 114
$ws.Cells.Item(152,3).Value = 0
$ws.Cells.Item(152,4).Value = 48
$ws.Cells.Item(152,5).Value = 60
$ws.Cells.Item(152,6).Value = 4
$ws.Cells.Item(152,8).Value = 6
$ws.Cells.Item(153,1).Value = 'Suazilandia'
$ws.Cells.Item(153,2).Value = 108
$ws.Cells.Item(153,3).Value = 2
$ws.Cells.Item(153,4).Value = 12
$ws.Cells.Item(153,5).Value = 95
$ws.Cells.Item(153,6).Value = 0
$ws.Cells.Item(153,8).Value = 1
$ws.Cells.Item(154,1).Value = 'Aruba'
$ws.Cells.Item(154,2).Value = 100
$ws.Cells.Item(154,4).Value = 81
$ws.Cells.Item(154,5).Value = 17
$ws.Cells.Item(154,6).Value = 4
$ws.Cells.Item(154,8).Value = 2
$ws.Cells.Item(155,1).Value = 'Monaco'
$ws.Cells.Item(155,2).Value = 95
$ws.Cells.Item(155,4).Value = 73
$ws.Cells.Item(155,5).Value = 18
$ws.Cells.Item(155,6).Value = 1
$ws.Cells.Item(155,8).Value = 4
$ws.Cells.Item(156,1).Value = 'Benin'
$ws.Cells.Item(156,2).Value = 90
$ws.Cells.Item(156,4).Value = 42
$ws.Cells.Item(156,5).Value = 46
$ws.Cells.Item(156,8).Value = 2
$ws.Cells.Item(157,1).Value = 'Haiti'
$ws.Cells.Item(157,4).Value = 10
$ws.Cells.Item(157,5).Value = 67
$ws.Cells.Item(157,8).Value = 8
$ws.Cells.Item(158,1).Value = 'Uganda'
$ws.Cells.Item(158,2).Value = 85
$ws.Cells.Item(158,4).Value = 52
$ws.Cells.Item(158,5).Value = 33
$ws.Cells.Item(158,6).Value = 0
$ws.Cells.Item(158,8).Value = 0
$ws.Cells.Item(159,1).Value = 'Guyana'
$ws.Cells.Item(159,3).Value = 0
$ws.Cells.Item(159,4).Value = 22
$ws.Cells.Item(159,5).Value = 51
$ws.Cells.Item(159,6).Value = 2
$ws.Cells.Item(159,8).Value = 9
$ws.Cells.Item(160,1).Value = 'Bahamas'
$ws.Cells.Item(160,3).Value = 1
$ws.Cells.Item(160,4).Value = 24
$ws.Cells.Item(160,5).Value = 47
$ws.Cells.Item(160,6).Value = 1
$ws.Cells.Item(160,8).Value = 11
$ws.Cells.Item(161,1).Value = 'Liechtenstein'
$ws.Cells.Item(161,2).Value = 82
$ws.Cells.Item(161,4).Value = 55
$ws.Cells.Item(161,5).Value = 26
$ws.Cells.Item(161,6).Value = 0
$ws.Cells.Item(161,8).Value = 1
$ws.Cells.Item(162,1).Value = 'Barbados'
$ws.Cells.Item(162,2).Value = 81
$ws.Cells.Item(162,4).Value = 44
$ws.Cells.Item(162,5).Value = 30
$ws.Cells.Item(162,6).Value = 4
$ws.Cells.Item(162,8).Value = 7
$ws.Cells.Item(163,1).Value = 'Mozambique'
$ws.Cells.Item(163,2).Value = 79
$ws.Cells.Item(163,4).Value = 18
$ws.Cells.Item(163,5).Value = 61
$ws.Cells.Item(163,6).Value = 0
$ws.Cells.Item(163,8).Value = 0
$ws.Cells.Item(164,1).Value = 'San Martin (Parte Holandesa)'
$ws.Cells.Item(164,2).Value = 76
$ws.Cells.Item(164,4).Value = 44
$ws.Cells.Item(164,5).Value = 19
$ws.Cells.Item(164,6).Value = 7
$ws.Cells.Item(164,8).Value = 13
$ws.Cells.Item(165,1).Value = 'Islas Caimanes'
$ws.Cells.Item(165,2).Value = 74
$ws.Cells.Item(165,4).Value = 10
$ws.Cells.Item(165,5).Value = 63
$ws.Cells.Item(165,6).Value = 3
$ws.Cells.Item(165,8).Value = 1
$ws.Cells.Item(178,1).Value = 'Angola'
$ws.Cells.Item(178,2).Value = 35
$ws.Cells.Item(178,3).Value = 5
$ws.Cells.Item(178,4).Value = 11
$ws.Cells.Item(178,5).Value = 22
$ws.Cells.Item(178,8).Value = 2
$ws.Cells.Item(179,1).Value = 'Zimbabue'
$ws.Cells.Item(179,2).Value = 34
$ws.Cells.Item(179,3).Value = 0
$ws.Cells.Item(179,4).Value = 5
$ws.Cells.Item(179,5).Value = 25
$ws.Cells.Item(179,8).Value = 4
$ws.Cells.Item(180,1).Value = 'Tayikistan'
$ws.Cells.Item(180,3).Value = 17
$ws.Cells.Item(180,5).Value = 32
$ws.Cells.Item(180,8).Value = 0
$ws.Cells.Item(181,1).Value = 'Guam'
$ws.Cells.Item(181,2).Value = 32
$ws.Cells.Item(181,4).Value = 0
$ws.Cells.Item(181,5).Value = 31
$ws.Cells.Item(181,8).Value = 1
$ws.Cells.Item(191,1).Value = 'Gambia'
$ws.Cells.Item(191,3).Value = 5
$ws.Cells.Item(191,4).Value = 9
$ws.Cells.Item(191,5).Value = 7
$ws.Cells.Item(191,8).Value = 1
$ws.Cells.Item(192,1).Value = 'Santa Lucia'
$ws.Cells.Item(192,2).Value = 17
$ws.Cells.Item(192,4).Value = 15
$ws.Cells.Item(192,5).Value = 2
$ws.Cells.Item(192,8).Value = 0
$ws.Cells.Item(193,1).Value = 'Santo Tome y Principe'
$ws.Cells.Item(193,4).Value = 4
$ws.Cells.Item(193,5).Value = 11
$ws.Cells.Item(193,8).Value = 1
$ws.Cells.Item(195,1).Value = 'San Vicente y las Granadinas'
$ws.Cells.Item(195,4).Value = 8
$ws.Cells.Item(195,5).Value = 8
$ws.Cells.Item(196,1).Value = 'Dominica'
$ws.Cells.Item(196,5).Value = 3
$ws.Cells.Item(196,8).Value = 0
$ws.Cells.Item(197,1).Value = 'Curazao'
$ws.Cells.Item(197,2).Value = 16
$ws.Cells.Item(197,4).Value = 13
$ws.Cells.Item(197,5).Value = 2
$ws.Cells.Item(197,8).Value = 1
$ws.Cells.Item(198,1).Value = 'San Cristobal y Nieves'
$ws.Cells.Item(198,2).Value = 15
$ws.Cells.Item(198,4).Value = 8
$ws.Cells.Item(198,5).Value = 7
$ws.Cells.Item(198,8).Value = 0
$ws.Cells.Item(199,1).Value = 'Nicaragua'
$ws.Cells.Item(199,2).Value = 14
$ws.Cells.Item(199,4).Value = 7
$ws.Cells.Item(199,5).Value = 4
$ws.Cells.Item(199,8).Value = 3
$ws.Cells.Item(200,1).Value = 'Islas Malvinas'
$ws.Cells.Item(200,2).Value = 13
$ws.Cells.Item(200,4).Value = 13
$ws.Cells.Item(200,5).Value = 0
$ws.Cells.Item(200,8).Value = 0
$ws.Cells.Item(201,1).Value = 'Islas Turcas y Caicos'
$ws.Cells.Item(201,4).Value = 5
$ws.Cells.Item(201,5).Value = 6
$ws.Cells.Item(207,1).Value = 'Yemen'
$ws.Cells.Item(207,3).Value = 3
$ws.Cells.Item(207,4).Value = 1
$ws.Cells.Item(207,5).Value = 7
$ws.Cells.Item(207,8).Value = 2
$ws.Cells.Item(208,1).Value = 'Surinam'
$ws.Cells.Item(208,2).Value = 10
$ws.Cells.Item(208,4).Value = 8
$ws.Cells.Item(208,5).Value = 1
$ws.Cells.Item(208,8).Value = 1
$ws.Cells.Item(209,1).Value = 'Papua Nueva Guinea'
$ws.Cells.Item(209,4).Value = 0
$ws.Cells.Item(209,5).Value = 8
$ws.Cells.Item(209,8).Value = 0
$ws.Cells.Item(210,1).Value = 'Mauritania'
$ws.Cells.Item(210,2).Value = 8
$ws.Cells.Item(210,4).Value = 6
$ws.Cells.Item(210,5).Value = 1
$ws.Cells.Item(210,8).Value = 1
$ws.Cells.Item(217,1).Value = 'Comoras'
$ws.Cells.Item(218,1).Value = 'San Pedro y Miquelon'
